$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("H1").Value = "Tempo Heuristica"
$ws.Range("I1").Value = "Tempo Total"

# Row 2 updates
$ws.Range("C2").Value = 1648
$ws.Range("F2").Value = "23 -> 22 -> 21 -> 67 -> 68 -> 70 -> 72 -> 10 -> 7 -> 4 -> 1 -> 2 -> 5 -> 8 -> 12 -> 16 -> 17 -> 20 -> 24 -> 23"
$ws.Range("G2").Value = 0.01498889923095703
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.01498889923095703

# Row 3 updates
$ws.Range("C3").Value = 2357
$ws.Range("F3").Value = "48 -> 42 -> 40 -> 39 -> 36 -> 35 -> 29 -> 13 -> 12 -> 8 -> 5 -> 4 -> 7 -> 10 -> 14 -> 15 -> 16 -> 17 -> 20 -> 24 -> 32 -> 41 -> 44 -> 46 -> 47 -> 48"
$ws.Range("G3").Value = 0.0273292064666748
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.0273292064666748

# Row 4 updates
$ws.Range("G4").Value = 0.02219200134277344
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.02219200134277344
